# 15-12-2016 -- 03:47 PM
# Add a new "View Records" test-case block (rows 35-43 and 47) to Sheet1,
# mirroring the layout/styling already used by the "Create New Lab" block
# (rows 1-15) and "Add Report" block (rows 18-32).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 35: section title bar (merged A35:C35), same look as A18:C18 ---
# Merge first (while still blank) so Excel has no pre-existing box border to
# redistribute, then paste the formatting over the now-merged range.
[void]$ws.Range("A35:C35").Merge()
$ws.Range("A18:C18").Copy()
$ws.Range("A35:C35").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A35").Value = "File Name : View Records."

# --- Row 36: "Mandatory Fields" label + value (B36:C36 merged) ---
[void]$ws.Range("B36:C36").Merge()
$ws.Range("A19").Copy()
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("B19:C19").Copy()
$ws.Range("B36:C36").PasteSpecial(-4122)
$ws.Range("A36").Value = "Mandatory Fields"
$ws.Range("B36").Value = "Phone,Passcode,MasterKey"

# --- Row 37: column headers (S.No / Test / Result / Return) ---
$ws.Range("A37").Value = "S.No"
$ws.Range("B37").Value = "Test"
$ws.Range("C37").Value = "Result / Return"

# --- Rows 38-43: the same six test-case rows as rows 6-11 ---
$ws.Range("C6").Copy()
$ws.Range("C38").PasteSpecial(-4122)

$ws.Range("A38").Value = 1
$ws.Range("B38").Value = "Empty any mandatory fields and form submit."
$ws.Range("C38").Value = "{""status"":false,""private key"":"""",""error"":""Lab Name or Person Name or Phone or passcode or master key is empty.""}"

$ws.Range("A39").Value = 2
$ws.Range("B39").Value = "Fill all mandatory fields but existing lab name."
$ws.Range("C39").Value = "{""status"":false,""private key"":"""",""error"":""Lab Name & Phone already exists.""}"

$ws.Range("A40").Value = 3
$ws.Range("B40").Value = "Fill all mandatory fields with new lab name and existing phone number and person name."
$ws.Range("C40").Value = "{""status"":true,""private key"":""847cfcec9cf78dd9cbf0dec71ac84851"",""error"":0}"

$ws.Range("A41").Value = 4
$ws.Range("B41").Value = "Fill all mandatory fields with existing lab name but new phone number and existing person name."
$ws.Range("C41").Value = "{""status"":true,""private key"":""b636fd739bb614bd17bb25210c68bed0"",""error"":0}"

$ws.Range("A42").Value = 5
$ws.Range("B42").Value = "Fill all mandatory fields with existing lab name and phone number but new person name."
$ws.Range("C42").Value = "{""status"":false,""private key"":"""",""error"":""Lab Name & Phone already exists.""}"

$ws.Range("A43").Value = 6
$ws.Range("B43").Value = "Fill all mandatory fields with new lab name and phone number but with wrong master key."
$ws.Range("C43").Value = "{""status"":false,""private key"":"""",""error"":""Master Key did not match.""}"

# --- Row 47: Result note (same text used on row 15) ---
$ws.Range("A47").Value = "Result"
$ws.Range("B47").Value = "If Lab Name or Lab Phone is new with correct masterkey data will be inserted."

[void]$ws.Range("B22").Select()

Write-Output "done"
